# Peak/Non-Peak pricing feature: update Pricing sheet rates, extend the
# Product sheet with Plan defaults, and rename AddOrder -> Order.

$wb = $excel.ActiveWorkbook

# --- Pricing sheet -------------------------------------------------------
$pricing = $wb.Worksheets.Item("Pricing")
$pricing.Range("C1").Value = "Web Data 13"
# Make room for the new "Time rate" style values by inserting a column
# before the existing WEEKDAY_BASED/SPECIAL_DAY/HOLIDAY/TIME_BASED block.
$pricing.Columns.Item(6).Insert()
$pricing.Range("F1").Value = "15"
$pricing.Range("F2").Value = "10"

# --- Product sheet ---------------------------------------------------------
$product = $wb.Worksheets.Item("Product")
$product.Range("I1").Value = "Plan"
$product.Range("J1").Value = "5"

# --- Rename AddOrder to Order ---------------------------------------------
$order = $wb.Worksheets.Item("AddOrder")
$order.Name = "Order"

# --- Make Product the active sheet ----------------------------------------
$product.Activate()
